$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old wrap-text formatting that column E inherited from the long
# "pesan" message column; we'll re-apply the correct (non-wrapping) look
# below, matching the rest of the sheet.
$ws.Columns.Item(5).ClearFormats()

# --- Column E header: "pesan" -> "file_pdf" -----------------------------
# Give E1 the same (non-wrapping) header formatting as the other header
# cells, then set the new label.
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "file_pdf"

# --- Column E body: long IF(...) message -> simple "<nisn>.pdf" formula --
# Give E2:E3 the same formatting as column B (vertical-center, no wrap,
# no horizontal centering) instead of the old wrapped "pesan" styling.
$ws.Range("B2:B3").Copy()
$ws.Range("E2:E3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E2").Formula = "=A2&"".pdf"""
$ws.Range("E3").Formula = "=A3&"".pdf"""

# The tall 45pt rows were only needed for the wrapped "pesan" text; restore
# the default row height now that the content is short again.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# Column E no longer needs to be extra wide to fit a long sentence.
$ws.Columns.Item(5).ColumnWidth = 16.14

# Match the saved selection state.
$ws.Range("E2:E3").Select()
